# 自动更新Excel文件 - 2025-11-01 23:12:10
# For every data row, the "剩余" (remaining) count in column E ticks down
# by one day. Once it would drop below 1, the cycle restarts at 10 and
# the "开始时间" (start date) in column F rolls forward by 10 days.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 99

for ($row = 2; $row -le $lastRow; $row++) {
    $eCell = $ws.Cells.Item($row, 5)
    $fCell = $ws.Cells.Item($row, 6)
    $remaining = $eCell.Value()
    $startRaw = $fCell.Value()

    if ($remaining -eq $null -or $remaining -eq "") {
        continue
    }

    # Skip rows whose start date isn't a well-formed yyyyMMdd value
    # (e.g. fat-fingered data entry) instead of corrupting them further.
    if ($startRaw -eq $null -or $startRaw -eq "" -or "$startRaw".Length -ne 8) {
        continue
    }

    $remaining = [int]$remaining

    if ($remaining -eq 1) {
        # Countdown hit zero: restock, reset the 10-day cycle and move
        # the start date forward by the cycle length.
        $startDate = [datetime]::ParseExact("$startRaw", "yyyyMMdd", $null)
        $newDate = $startDate.AddDays(10)
        $fCell.Value = [int]$newDate.ToString("yyyyMMdd")
        $eCell.Value = 10
    }
    else {
        $eCell.Value = $remaining - 1
    }
}
